$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 461.05884
$ws.Range("I107").Value = 353.54544
$ws.Range("J107").Value = 658.1667
$ws.Range("K107").Value = 353.54544
$ws.Range("L107").Value = 658.1667
$ws.Range("M107").Value = 1566.45456
$ws.Range("N107").Value = -4498.1667
$ws.Range("H111").Value = 443.75
$ws.Range("I111").Value = 425
$ws.Range("J111").Value = 462.5
$ws.Range("K111").Value = 1275
$ws.Range("L111").Value = 1387.5
$ws.Range("M111").Value = 1792
$ws.Range("N111").Value = -7521.5
$ws.Range("H132").Value = 6949018
$ws.Range("I132").Value = 7579435
$ws.Range("J132").Value = 14433
$ws.Range("K132").Value = 22738305
$ws.Range("L132").Value = 43299
$ws.Range("M132").Value = -22735775
$ws.Range("N132").Value = -48359
$ws.Range("H135").Value = 911.3953
$ws.Range("I135").Value = 873.5714
$ws.Range("J135").Value = 2500
$ws.Range("K135").Value = 7862.1426
$ws.Range("L135").Value = 22500
$ws.Range("M135").Value = -5327.1426
$ws.Range("N135").Value = -27570
$ws.Range("H138").Value = 1571.386
$ws.Range("I138").Value = 965.97675
$ws.Range("J138").Value = 3430.8572
$ws.Range("K138").Value = 2897.93025
$ws.Range("L138").Value = 10292.5716
$ws.Range("M138").Value = 2242.06975
$ws.Range("N138").Value = -20572.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 904.17
$ws.Range("I32").Value = 813.6667
$ws.Range("J32").Value = 1718.7
$ws.Range("K32").Value = 813.6667
$ws.Range("L32").Value = 1718.7
$ws.Range("M32").Value = -526.6667
$ws.Range("N32").Value = -2292.7
$ws.Range("H61").Value = 1383.3667
$ws.Range("I61").Value = 1415.8948
$ws.Range("J61").Value = 1327.1818
$ws.Range("K61").Value = 1415.8948
$ws.Range("L61").Value = 1327.1818
$ws.Range("M61").Value = -1203.8948
$ws.Range("N61").Value = -1751.1818
$ws.Range("H136").Value = 1383.3667
$ws.Range("I136").Value = 1415.8948
$ws.Range("J136").Value = 1327.1818
$ws.Range("K136").Value = 4247.6844
$ws.Range("L136").Value = 3981.5454
$ws.Range("M136").Value = -1697.6844
$ws.Range("N136").Value = -9081.545399999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1064.5151
$ws.Range("I80").Value = 824.1
$ws.Range("J80").Value = 1169.0435
$ws.Range("K80").Value = 824.1
$ws.Range("L80").Value = 1169.0435
$ws.Range("M80").Value = 173.9
$ws.Range("N80").Value = -3165.0435
$ws.Range("H83").Value = 1064.5151
$ws.Range("I83").Value = 824.1
$ws.Range("J83").Value = 1169.0435
$ws.Range("K83").Value = 4120.5
$ws.Range("L83").Value = 5845.2175
$ws.Range("M83").Value = 871.5
$ws.Range("N83").Value = -15829.2175
$ws.Range("H134").Value = 1986732.8
$ws.Range("I134").Value = 678.7368
$ws.Range("J134").Value = 6179513.5
$ws.Range("K134").Value = 2036.2104
$ws.Range("L134").Value = 18538540.5
$ws.Range("M134").Value = 498.7896000000001
$ws.Range("N134").Value = -18543610.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 19232122
$ws.Range("I58").Value = 25642556
$ws.Range("J58").Value = 819.53845
$ws.Range("K58").Value = 25642556
$ws.Range("L58").Value = 819.53845
$ws.Range("M58").Value = -25642353
$ws.Range("N58").Value = -1225.53845
$ws.Range("H136").Value = 19232122
$ws.Range("I136").Value = 25642556
$ws.Range("J136").Value = 819.53845
$ws.Range("K136").Value = 76927668
$ws.Range("L136").Value = 2458.61535
$ws.Range("M136").Value = -76925118
$ws.Range("N136").Value = -7558.61535

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 33336810
$ws.Range("I5").Value = 35088096
$ws.Range("J5").Value = 31257156
$ws.Range("K5").Value = 105264288
$ws.Range("L5").Value = 93771468
$ws.Range("M5").Value = -105264176
$ws.Range("N5").Value = -93771692
$ws.Range("H122").Value = 8625018
$ws.Range("I122").Value = 35714564
$ws.Range("J122").Value = 5617.727
$ws.Range("K122").Value = 321431076
$ws.Range("L122").Value = 50559.543
$ws.Range("M122").Value = -321428626
$ws.Range("N122").Value = -55459.543
$ws.Range("H131").Value = 792.58
$ws.Range("I131").Value = 529.75
$ws.Range("J131").Value = 815.43475
$ws.Range("K131").Value = 1589.25
$ws.Range("L131").Value = 2446.30425
$ws.Range("M131").Value = 3450.75
$ws.Range("N131").Value = -12526.30425
$ws.Range("H135").Value = 33336810
$ws.Range("I135").Value = 35088096
$ws.Range("J135").Value = 31257156
$ws.Range("K135").Value = 315792864
$ws.Range("L135").Value = 281314404
$ws.Range("M135").Value = -315790329
$ws.Range("N135").Value = -281319474

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7695792.5
$ws.Range("I80").Value = 4950.8335
$ws.Range("J80").Value = 14287943
$ws.Range("K80").Value = 4950.8335
$ws.Range("L80").Value = 14287943
$ws.Range("M80").Value = -3952.8335
$ws.Range("N80").Value = -14289939
$ws.Range("H83").Value = 7695792.5
$ws.Range("I83").Value = 4950.8335
$ws.Range("J83").Value = 14287943
$ws.Range("K83").Value = 24754.1675
$ws.Range("L83").Value = 71439715
$ws.Range("M83").Value = -19762.1675
$ws.Range("N83").Value = -71449699
$ws.Range("H97").Value = 876.1070999999999
$ws.Range("I97").Value = 651.875
$ws.Range("J97").Value = 2221.5
$ws.Range("K97").Value = 651.875
$ws.Range("L97").Value = 2221.5
$ws.Range("M97").Value = -155.875
$ws.Range("N97").Value = -3213.5
$ws.Range("H102").Value = 2536.2856
$ws.Range("I102").Value = 2149.6
$ws.Range("J102").Value = 3503
$ws.Range("K102").Value = 2149.6
$ws.Range("L102").Value = 3503
$ws.Range("M102").Value = -527.5999999999999
$ws.Range("N102").Value = -6747
$ws.Range("H132").Value = 6036.8887
$ws.Range("I132").Value = 3676.7097
$ws.Range("J132").Value = 20670
$ws.Range("K132").Value = 11030.1291
$ws.Range("L132").Value = 62010
$ws.Range("M132").Value = -8500.1291
$ws.Range("N132").Value = -67070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2819.8948
$ws.Range("I136").Value = 3277.0356
$ws.Range("J136").Value = 1539.9
$ws.Range("K136").Value = 9831.106800000001
$ws.Range("L136").Value = 4619.700000000001
$ws.Range("M136").Value = -7281.106800000001
$ws.Range("N136").Value = -9719.700000000001
